# Weekly update: insert a new "Damasco" price row at the top of the
# dated data (row 10), pushing the previously-existing rows 10-13 down
# to rows 11-14. The new row carries this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 10:13 down to 11:14 and make room for the new entry.
$ws.Rows("10:10").Insert()

# Populate the new row 10 with this week's observation.
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").Value = 44917
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100103
$ws.Range("H10").Value = "Frutos de hueso (carozo)"
$ws.Range("I10").Value = 100103003
$ws.Range("J10").Value = "Damasco"
$ws.Range("K10").Value = "Castle Brite"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 23000
$ws.Range("P10").Value = 21800
$ws.Range("Q10").Value = "`$/caja 18 kilos"
$ws.Range("R10").Value = "Región de Coquimbo"
$ws.Range("S10").Value = 1211
$ws.Range("T10").Value = 18
